# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row 2 for the 2022-Q3 summary,
#    pushing the existing 2022-Q2 / 2022-Q1 rows down by one.
# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e. before
#    the existing "2022-Q2" sheet) containing the Q3 per-fund holdings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert new summary row for 2022-Q3
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-use the formatting already on the (now shifted) 2022-Q2 row for the A
# column "index" cell, matching the s="2" style used throughout column A.
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# The A-column is a 0-based row index; renumber the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Re-fetch by name *after* the Add() shifted sheet positions, so this
# resolves to the actual "2022-Q2" sheet rather than whatever now sits in
# its old slot.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Match header/first-column formatting used on the other quarter sheets.
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q2Sheet.Range("A2").Copy($newSheet.Range("A2"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# The fund-code / size / position columns are stored as text (leading
# zeros in fund codes must survive), matching the other quarter sheets.
$newSheet.Range("B2:G2").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "001614"
$newSheet.Range("C2").Value = "东方区域发展混合"
$newSheet.Range("D2").Value = "0.22"
$newSheet.Range("E2").Value = "99.06"
$newSheet.Range("F2").Value = "5.71"
$newSheet.Range("G2").Value = "0.0126"
$newSheet.Range("H2").Value = 10
